$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.219803940679071
$ws.Range("D2").Value = 0.00136281687266937
$ws.Range("E2").Value = 0.4336984297763706
$ws.Range("F2").Value = 0.7160575709075232
$ws.Range("G2").Value = 0.002324725457951165
$ws.Range("I2").Value = 3.961193417866525
$ws.Range("O2").Value = 2.355210718842727
$ws.Range("B3").Value = 1.069997499266549
$ws.Range("D3").Value = 0.001187382742090293
$ws.Range("E3").Value = 0.3778594790115903
$ws.Range("F3").Value = 0.6567163053088905
$ws.Range("G3").Value = 0.002329461130050568
$ws.Range("I3").Value = 3.481483530427994
$ws.Range("O3").Value = 2.156718067768907
$ws.Range("B4").Value = 0.9775436500530645
$ws.Range("D4").Value = 0.001080080274057593
$ws.Range("E4").Value = 0.3437100912600357
$ws.Range("F4").Value = 0.6207120828670867
$ws.Range("G4").Value = 0.002332519531335734
$ws.Range("I4").Value = 3.186307868167489
$ws.Range("O4").Value = 2.036265997455075
$ws.Range("B5").Value = 0.9397517899408854
$ws.Range("D5").Value = 0.001036432189167158
$ws.Range("E5").Value = 0.329825076174572
$ws.Range("F5").Value = 0.6061468742205847
$ws.Range("G5").Value = 0.002333803886885216
$ws.Range("I5").Value = 3.065876933603903
$ws.Range("O5").Value = 1.987532592803802
$ws.Range("B6").Value = 0.9334695415156489
$ws.Range("D6").Value = 0.001029188458627317
$ws.Range("E6").Value = 0.3275212721868144
$ws.Range("F6").Value = 0.6037347334826251
$ws.Range("G6").Value = 0.002334019454337199
$ws.Range("I6").Value = 3.04587117873993
$ws.Range("O6").Value = 1.979461519004303
$ws.Range("B7").Value = 0.9770344427063833
$ws.Range("D7").Value = 0.001079491335964633
$ws.Range("E7").Value = 0.3435227109891912
$ws.Range("F7").Value = 0.6205152211821883
$ws.Range("G7").Value = 0.002332536698358234
$ws.Range("I7").Value = 3.184684259858585
$ws.Range("O7").Value = 2.035607345351195
$ws.Range("B8").Value = 1.168250105383493
$ws.Range("D8").Value = 0.001302225393729373
$ws.Range("E8").Value = 0.4144150822294819
$ws.Range("F8").Value = 0.6955060775010651
$ws.Range("G8").Value = 0.002326327130626193
$ws.Range("I8").Value = 3.795928923952232
$ws.Range("O8").Value = 2.286471316272696
$ws.Range("B9").Value = 1.539389405566169
$ws.Range("D9").Value = 0.001743495443395915
$ws.Range("E9").Value = 0.5546583365363915
$ws.Range("F9").Value = 0.8460728912957052
$ws.Range("G9").Value = 0.002315339238254839
$ws.Range("I9").Value = 4.989014804984947
$ws.Range("O9").Value = 2.79001036887513
$ws.Range("B10").Value = 1.809644099351658
$ws.Range("D10").Value = 0.002072145766488731
$ws.Range("E10").Value = 0.658652049917535
$ws.Range("F10").Value = 0.9589651656520459
$ws.Range("G10").Value = 0.002307982165642218
$ws.Range("I10").Value = 5.861568061238756
$ws.Range("O10").Value = 3.167490004210322
$ws.Range("B11").Value = 1.932050532542917
$ws.Range("D11").Value = 0.002223006803760086
$ws.Range("E11").Value = 0.7062161564137739
$ws.Range("F11").Value = 1.010845266279489
$ws.Range("G11").Value = 0.002304788697644411
$ws.Range("I11").Value = 6.25753114434093
$ws.Range("O11").Value = 3.340953389518518
$ws.Range("B12").Value = 1.978324344473322
$ws.Range("D12").Value = 0.002280360919026947
$ws.Range("E12").Value = 0.7242680764979781
$ws.Range("F12").Value = 1.03056856590527
$ws.Range("G12").Value = 0.00230360130780261
$ws.Range("I12").Value = 6.407322560669854
$ws.Range("O12").Value = 3.406898250946483
$ws.Range("B13").Value = 1.968361989860625
$ws.Range("D13").Value = 0.002267998111934233
$ws.Range("E13").Value = 0.720378418342051
$ws.Range("F13").Value = 1.026317325067765
$ws.Range("G13").Value = 0.00230385606152078
$ws.Range("I13").Value = 6.375069214559744
$ws.Range("O13").Value = 3.392684256368284
$ws.Range("B14").Value = 1.935859098657602
$ws.Range("D14").Value = 0.00222772064371668
$ws.Range("E14").Value = 0.7077004683360144
$ws.Range("F14").Value = 1.012466350817789
$ws.Range("G14").Value = 0.002304690571907265
$ws.Range("I14").Value = 6.269857673435979
$ws.Range("O14").Value = 3.346373501059304
$ws.Range("B15").Value = 1.915939811075305
$ws.Range("D15").Value = 0.002203079946879072
$ws.Range("E15").Value = 0.6999402271471666
$ws.Range("F15").Value = 1.003992371776604
$ws.Range("G15").Value = 0.002305204582768017
$ws.Range("I15").Value = 6.205392518134659
$ws.Range("O15").Value = 3.318040639829917
$ws.Range("B16").Value = 1.801633628443199
$ws.Range("D16").Value = 0.002062316590961899
$ws.Range("E16").Value = 0.655549089093725
$ws.Range("F16").Value = 0.9555854069586189
$ws.Range("G16").Value = 0.002308193939223566
$ws.Range("I16").Value = 5.835670441680065
$ws.Range("O16").Value = 3.15618949106215
$ws.Range("B17").Value = 1.731372246120657
$ws.Range("D17").Value = 0.001976332213386556
$ws.Range("E17").Value = 0.6283848123191405
$ws.Range("F17").Value = 0.9260251098173171
$ws.Range("G17").Value = 0.002310066979291825
$ws.Range("I17").Value = 5.608601612490133
$ws.Range("O17").Value = 3.057351167739682
$ws.Range("B18").Value = 1.690909632147168
$ws.Range("D18").Value = 0.001927002019733948
$ws.Range("E18").Value = 0.6127846016899383
$ws.Range("F18").Value = 0.9090720897889923
$ws.Range("G18").Value = 0.00231115873986465
$ws.Range("I18").Value = 5.477907426105048
$ws.Range("O18").Value = 3.000665896692624
$ws.Range("B19").Value = 1.677201143001923
$ws.Range("D19").Value = 0.001910320321501402
$ws.Range("E19").Value = 0.6075066487272522
$ws.Range("F19").Value = 0.903340493389095
$ws.Range("G19").Value = 0.002311530874994851
$ws.Range("I19").Value = 5.433641513037287
$ws.Range("O19").Value = 2.981501180670932
$ws.Range("B20").Value = 1.738856893519369
$ws.Range("D20").Value = 0.001985472159233836
$ws.Range("E20").Value = 0.6312739907922804
$ws.Range("F20").Value = 0.9291667394193155
$ws.Range("G20").Value = 0.002309866097720437
$ws.Range("I20").Value = 5.632782903627003
$ws.Range("O20").Value = 3.067855660969485
$ws.Range("B21").Value = 1.945408141630878
$ws.Range("D21").Value = 0.002239544716820419
$ws.Range("E21").Value = 0.7114231643078455
$ws.Range("F21").Value = 1.016532600446681
$ws.Range("G21").Value = 0.002304444862530205
$ws.Range("I21").Value = 6.300765052729901
$ws.Range("O21").Value = 3.359969033151685
$ws.Range("B22").Value = 2.079940855283837
$ws.Range("D22").Value = 0.002406930069309254
$ws.Range("E22").Value = 0.7640425772889614
$ws.Range("F22").Value = 1.074083505283141
$ws.Range("G22").Value = 0.002301029401395692
$ws.Range("I22").Value = 6.736443306619549
$ws.Range("O22").Value = 3.552389396738079
$ws.Range("B23").Value = 2.008181005424092
$ws.Range("D23").Value = 0.00231746072000405
$ws.Range("E23").Value = 0.7359357237643849
$ws.Range("F23").Value = 1.043325465387227
$ws.Range("G23").Value = 0.002302840664474282
$ws.Range("I23").Value = 6.503998686294949
$ws.Range("O23").Value = 3.449550766766095
$ws.Range("B24").Value = 1.735473297350381
$ws.Range("D24").Value = 0.001981339670157567
$ws.Range("E24").Value = 0.6299677408469222
$ws.Range("F24").Value = 0.927746279284861
$ws.Range("G24").Value = 0.002309956869657329
$ws.Range("I24").Value = 5.62185100597145
$ws.Range("O24").Value = 3.06310614965065
$ws.Range("B25").Value = 1.439406031740532
$ws.Range("D25").Value = 0.00162346534227531
$ws.Range("E25").Value = 0.5165662124021679
$ws.Range("F25").Value = 0.804951376280016
$ws.Range("G25").Value = 0.002318185404700003
$ws.Range("I25").Value = 4.666914883151151
$ws.Range("O25").Value = 2.652501122830984
